$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of COVID overview data for the UK (2022-01-23 .. 2022-01-30),
# appended after the existing last row (529).
$rows = @(
    @{ Row = 530; Date = "2022-01-23"; CumCases = 15859288; NewCases = 74799; NewDeaths = 75;  CumDeaths = 153862 },
    @{ Row = 531; Date = "2022-01-24"; CumCases = 15953685; NewCases = 88447; NewDeaths = 56;  CumDeaths = 153916 },
    @{ Row = 532; Date = "2022-01-25"; CumCases = 16047716; NewCases = 94326; NewDeaths = 439; CumDeaths = 154356 },
    @{ Row = 533; Date = "2022-01-26"; CumCases = 16149319; NewCases = 102292; NewDeaths = 346; CumDeaths = 154702 },
    @{ Row = 534; Date = "2022-01-27"; CumCases = 16245474; NewCases = 96871; NewDeaths = 338; CumDeaths = 155040 },
    @{ Row = 535; Date = "2022-01-28"; CumCases = 16333980; NewCases = 89176; NewDeaths = 277; CumDeaths = 155317 },
    @{ Row = 536; Date = "2022-01-29"; CumCases = 16406123; NewCases = 72727; NewDeaths = 296; CumDeaths = 155613 },
    @{ Row = 537; Date = "2022-01-30"; CumCases = 16468522; NewCases = 62399; NewDeaths = 85;  CumDeaths = 155698 }
)

# Make sure the date column (A) is formatted as text before assigning the
# date-like strings, otherwise Excel auto-converts "2022-01-23" into a
# date serial number instead of keeping it as literal text (matching the
# existing rows above, which all store the date as plain text).
$ws.Range("A530:A537").NumberFormat = "@"

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "overview"
    $ws.Cells.Item($row, 3).Value = "K02000001"
    $ws.Cells.Item($row, 4).Value = "United Kingdom"
    $ws.Cells.Item($row, 5).Value = $r.CumCases
    $ws.Cells.Item($row, 6).Value = $r.NewCases
    $ws.Cells.Item($row, 7).Value = $r.NewDeaths
    $ws.Cells.Item($row, 8).Value = $r.CumDeaths
}
